$wb = $excel.ActiveWorkbook

# Both the "展览" sheet and the "全部类型" sheet mirror the same event data,
# so the "想去人数" (interest count) update needs to land on both of them.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 299
    $ws.Range("F3").Value = 167
}
